# Corrige les noms des paroisses : remplace le trait d'union par un tiret
# demi-cadratin entoure d'espaces ("Pully-Paudex" -> "Pully – Paudex", etc.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "Pully – Paudex"
$ws.Range("E3").Value = "Belmont – Lutry"
$ws.Range("E5").Value = "Savigny – Forel"

$ws.Range("E5").Select()
